$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = '97.294.74'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '3.729.17'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '238.92'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").Value = '1.89'
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("D7").Value = '662.48'
$ws.Range("E7").Value = '  +1.08%  '
$ws.Range("D8").Value = '0.437'
$ws.Range("E8").Value = '  +3.18%  '
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D11").Value = '3.729.28'
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").Value = '0.0000324'
$ws.Range("E12").Value = '  +20.34%  '
$ws.Range("D13").Value = '45.22'
$ws.Range("E14").Value = '  +1.69%  '
$ws.Range("D15").Value = '6.95'
$ws.Range("E15").Value = '  +1.93%  '
$ws.Range("D16").Value = '4.426.06'
$ws.Range("D17").Value = '97.016.93'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '9.10'
$ws.Range("E18").Value = '  +2.60%  '
$ws.Range("D19").Value = '3.722.81'
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = '13.13'
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("D21").Value = '18.85'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '0.507'
$ws.Range("E22").Value = '  -3.75%  '
$ws.Range("D23").Value = '527.80'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '3.51'
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("D25").Value = '0.0000230'
$ws.Range("E25").Value = '  +12.56%  '
$ws.Range("E26").Value = '  -2.43%  '
$ws.Range("D27").Value = '109.26'
$ws.Range("E27").Value = '  +6.96%  '
$ws.Range("D28").Value = '0.193'
$ws.Range("E28").Value = '  +14.64%  '
$ws.Range("D29").Value = '13.69'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").Value = '3.933.66'
$ws.Range("E30").Value = '  +1.55%  '
$ws.Range("D31").Value = '12.96'
$ws.Range("E31").Value = '  +3.81%  '
$ws.Range("D32").Value = '3.07'
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '0.193'
$ws.Range("E34").Value = '  +4.32%  '
$ws.Range("D35").Value = '1.85'
$ws.Range("E35").Value = '  -4.15%  '
$ws.Range("D36").Value = '32.64'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").Value = '647.48'
$ws.Range("E38").Value = '  -4.11%  '
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("D40").Value = '8.85'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '0.168'
$ws.Range("E42").Value = '  +4.26%  '
$ws.Range("D46").Value = '0.982'
$ws.Range("E46").Value = '  +2.10%  '
$ws.Range("D47").Value = '0.483'
$ws.Range("E47").Value = '  +10.25%  '
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").Value = '2.40'
$ws.Range("E49").Value = '  +3.26%  '

# --- Row reordering (B/C identity + D/E values) for rows 43-45 and 50-51 ---
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '6.85'
$ws.Range("E43").Value = '  +4.58%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = '2.05'
$ws.Range("E44").Value = '  +2.90%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '40.82'
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '23.63'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '8.76'
$ws.Range("E51").Value = '  +1.23%  '
